# The deck originally had two extra slides near the end -- "Conclusions"
# (slide 17, sldId 277) and "References" (slide 18, sldId 278) -- that were
# removed, leaving "Thank you" (previously slide 19, sldId 279) as the new
# final slide.
#
# Deleting slide 17 first shifts "References" into position 17, so deleting
# position 17 a second time removes it too and "Thank you" naturally settles
# into position 17.

$p = $ppt.ActivePresentation

$p.Slides.Item(17).Delete()
$p.Slides.Item(17).Delete()
